$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for specific rows per repulled data
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = 3
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -8
$ws.Range("F10").Value = 4
